$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-coerced to a number by Excel
# (single decimal point, no letters) need NumberFormat forced to Text first
# so they remain text cells, matching the source data feed which always
# writes these as strings.
$ws.Range("D2").Value = "59.344.46"
$ws.Range("E2").Value = "  +2.79%  "
$ws.Range("D3").Value = "3.189.05"
$ws.Range("E3").Value = "  +2.02%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "532.61"
$ws.Range("E5").Value = "  +0.02%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.38"
$ws.Range("E6").Value = "  +2.70%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.525"
$ws.Range("E8").Value = "  +10.44%  "
$ws.Range("E9").Value = "  -0.33%  "
$ws.Range("E10").Value = "  +6.19%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.112"
$ws.Range("E11").Value = "  +4.79%  "
$ws.Range("B12").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C12").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D12").Value = "3.736.45"
$ws.Range("E12").Value = "  +2.08%  "
$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.140"
$ws.Range("E13").Value = "  +1.88%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.03"
$ws.Range("E14").Value = "  +0.47%  "
$ws.Range("E15").Value = "  +4.75%  "
$ws.Range("D16").Value = "59.401.07"
$ws.Range("E16").Value = "  +2.69%  "
$ws.Range("B17").Value = "Polkadot"
$ws.Range("C17").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.26"
$ws.Range("E17").Value = "  +2.89%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.171.45"
$ws.Range("E18").Value = "  +1.28%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.08"
$ws.Range("E19").Value = "  +2.87%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.22"
$ws.Range("E20").Value = "  +1.84%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "375.66"
$ws.Range("E21").Value = "  +2.14%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.999"
$ws.Range("E22").Value = "  +0.16%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.534"
$ws.Range("E23").Value = "  +5.50%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "69.87"
$ws.Range("E25").Value = "  +0.34%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.996"
$ws.Range("E26").Value = "  -0.41%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.49"
$ws.Range("E27").Value = "  +16.21%  "
$ws.Range("D28").Value = "0.0₃0878"
$ws.Range("E28").Value = "  +1.64%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "22.43"
$ws.Range("E29").Value = "  +4.68%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.90"
$ws.Range("E30").Value = "  +1.24%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.11"
$ws.Range("E31").Value = "  +0.72%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.25"
$ws.Range("E32").Value = "  +1.99%  "
$ws.Range("E33").Value = "  -0.06%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.37"
$ws.Range("E34").Value = "  +4.76%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "157.22"
$ws.Range("E35").Value = "  -1.60%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.35"
$ws.Range("E36").Value = "  +3.99%  "
$ws.Range("B37").Value = "EnergySwap"
$ws.Range("C37").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "25.66"
$ws.Range("E37").Value = "  +1.02%  "
$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0715"
$ws.Range("E38").Value = "  +6.18%  "
$ws.Range("D39").Value = "2.718.26"
$ws.Range("E39").Value = "  +7.28%  "
$ws.Range("E40").Value = "  +1.96%  "
$ws.Range("E41").Value = "  +4.90%  "
$ws.Range("B42").Value = "Mantle"
$ws.Range("C42").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.728"
$ws.Range("E42").Value = "  +4.17%  "
$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0294"
$ws.Range("E43").Value = "  +8.91%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "39.27"
$ws.Range("E44").Value = "  +4.24%  "
$ws.Range("E45").Value = "  +0.08%  "
$ws.Range("D46").Value = "3.233.25"
$ws.Range("E46").Value = "  +2.03%  "
$ws.Range("B47").Value = "ONDO"
$ws.Range("C47").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.992"
$ws.Range("E47").Value = "  +1.71%  "
$ws.Range("B48").Value = "Stellar"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.102"
$ws.Range("E48").Value = "  +11.88%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.24"
$ws.Range("E49").Value = "  +1.83%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "20.47"
$ws.Range("E50").Value = "  +3.51%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.761"
$ws.Range("E51").Value = "  +2.42%  "
